$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC01_CCDI_phs002371_Sex-Male.xlsx — "Fixed CCDI queries and ICDC index issue"
# The "Participants" tab query (cell B2) is rewritten: it drops the
# "Ethnicity" column and renames "Alternate ID" to "Synonym Participant ID".
$newQuery = @'
SELECT
    p.participant_id AS "Participant ID",
    st.study_id AS "Study ID",
    COALESCE(p.sex_at_birth, '') AS "Sex",
    COALESCE(p.race, '') AS "Race",
    COALESCE(CAST(p.alternate_participant_id AS INT), '') AS "Synonym Participant ID"
FROM 
    df_participant p
JOIN 
    df_study st ON p."study.id" = st.id
LEFT JOIN 
    df_sample smp ON smp."participant.id" = p.participant_id
LEFT JOIN 
    df_diagnosis diag ON diag."participant.id" = p.participant_id
LEFT JOIN 
    df_publication pub ON pub."study.id" = st.study_id
LEFT JOIN 
    df_sequencing_file sqf ON sqf."sample.id" = smp.sample_id
WHERE 
    st.study_id = 'phs002371' AND p.sex_at_birth = 'Male'
ORDER BY 
    p.participant_id ASC
LIMIT 100;
'@

$ws.Range("B2").Value2 = $newQuery

# Move the active selection to C2, matching the saved workbook state.
$ws.Range("C2").Select()
